$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D:E columns to text so numeric-looking strings (e.g. "0.536") are not
# auto-converted to numbers by Excel; style is restored to Normal afterwards.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '58.839.11'
$ws.Range('E2').Value = '  +1.48%  '
$ws.Range('D3').Value = '3.159.09'
$ws.Range('E3').Value = '  +1.15%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '533.13'
$ws.Range('E5').Value = '  +0.58%  '
$ws.Range('D6').Value = '140.38'
$ws.Range('E6').Value = '  +1.29%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '0.536'
$ws.Range('E8').Value = '  +16.26%  '
$ws.Range('D9').Value = '7.33'
$ws.Range('E9').Value = '  +0.57%  '
$ws.Range('D10').Value = '0.432'
$ws.Range('E10').Value = '  +5.94%  '
$ws.Range('D11').Value = '0.111'
$ws.Range('E11').Value = '  +3.58%  '
$ws.Range('E12').Value = '  +2.70%  '
$ws.Range('D13').Value = '3.701.25'
$ws.Range('E13').Value = '  +1.08%  '
$ws.Range('D14').Value = '26.03'
$ws.Range('E14').Value = '  +2.01%  '
$ws.Range('D15').Value = '0.0000172'
$ws.Range('E15').Value = '  +5.68%  '
$ws.Range('D16').Value = '58.859.33'
$ws.Range('E16').Value = '  +1.45%  '
$ws.Range('D17').Value = '6.26'
$ws.Range('E17').Value = '  +4.68%  '
$ws.Range('D18').Value = '3.157.49'
$ws.Range('E18').Value = '  +0.89%  '
$ws.Range('D19').Value = '13.03'
$ws.Range('E19').Value = '  +2.88%  '
$ws.Range('D20').Value = '8.22'
$ws.Range('E20').Value = '  +2.00%  '
$ws.Range('D21').Value = '371.91'
$ws.Range('E21').Value = '  +5.61%  '
$ws.Range('E22').Value = '  +2.01%  '
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').Value = '70.07'
$ws.Range('E24').Value = '  +1.97%  '
$ws.Range('D25').Value = '0.520'
$ws.Range('E25').Value = '  +3.30%  '
$ws.Range('E26').Value = '  +0.40%  '
$ws.Range('E27').Value = '  -0.24%  '
$ws.Range('D28').Value = '8.17'
$ws.Range('E28').Value = '  +12.30%  '
$ws.Range('D29').Value = '0.0₃0867'
$ws.Range('E29').Value = '  -1.94%  '
$ws.Range('D30').Value = '1.89'
$ws.Range('E30').Value = '  +1.71%  '
$ws.Range('D31').Value = '6.13'
$ws.Range('E31').Value = '  +0.46%  '
$ws.Range('D32').Value = '22.09'
$ws.Range('E32').Value = '  +3.75%  '
$ws.Range('D33').Value = '5.21'
$ws.Range('E33').Value = '  +4.43%  '
$ws.Range('D34').Value = '1.17'
$ws.Range('E34').Value = '  +1.16%  '
$ws.Range('D35').Value = '159.20'
$ws.Range('E35').Value = '  +0.51%  '
$ws.Range('D36').Value = '6.29'
$ws.Range('E36').Value = '  +3.85%  '
$ws.Range('E37').Value = '  +7.45%  '
$ws.Range('D38').Value = '25.26'
$ws.Range('E38').Value = '  -3.11%  '
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('D40').Value = '2.641.41'
$ws.Range('E40').Value = '  +10.44%  '
$ws.Range('D41').Value = '0.0684'
$ws.Range('E41').Value = '  +2.25%  '
$ws.Range('D42').Value = '4.25'
$ws.Range('E42').Value = '  +6.35%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '0.0288'
$ws.Range('E43').Value = '  +7.87%  '
$ws.Range('D44').Value = '0.713'
$ws.Range('E44').Value = '  +2.19%  '
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').Value = '38.80'
$ws.Range('E45').Value = '  +3.44%  '
$ws.Range('E46').Value = '  -0.02%  '
$ws.Range('D47').Value = '3.198.28'
$ws.Range('E47').Value = '  +1.04%  '
$ws.Range('D48').Value = '0.104'
$ws.Range('E48').Value = '  +14.47%  '
$ws.Range('D49').Value = '0.986'
$ws.Range('E49').Value = '  +1.33%  '
$ws.Range('D50').Value = '6.21'
$ws.Range('E50').Value = '  +2.87%  '
$ws.Range('D51').Value = '20.31'
$ws.Range('E51').Value = '  +2.52%  '

$priceRange.Style = "Normal"
